$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StudentLogin")
$ws.Activate()

# --- Capture the original style objects so we can restore them after the
# hyperlink operations below (adding a hyperlink forces Excel to apply a
# possibly-new "Hyperlink" style index to the cell). Rows 6/7 col A already
# use the built-in Hyperlink look (style index 1), so remember it.
$styleA6 = $ws.Range("A6").Style
$styleA7 = $ws.Range("A7").Style

# --- Swap the contents of row 6 and row 7 (columns A:C) ---
$a6 = $ws.Range("A6").Value()
$b6 = $ws.Range("B6").Value()
$c6 = $ws.Range("C6").Value()
$a7 = $ws.Range("A7").Value()
$b7 = $ws.Range("B7").Value()
$c7 = $ws.Range("C7").Value()

$ws.Range("A6").Value = $a7
$ws.Range("B6").Value = $b7
$ws.Range("C6").Value = $c7

$ws.Range("A7").Value = $a6
$ws.Range("B7").Value = $b6
$ws.Range("C7").Value = $c6

# --- Remove the two existing hyperlinks that lived on A6 and A7.
# Deleting from the Hyperlinks collection while enumerating it shifts the
# remaining items, so re-scan from scratch after every removal. ---
function Remove-HyperlinkAtAddress($sheet, $targetAddress) {
    $found = $true
    while ($found) {
        $found = $false
        foreach ($hl in $sheet.Hyperlinks) {
            if ($hl.Range.Address() -eq $targetAddress) {
                $hl.Delete()
                $found = $true
                break
            }
        }
    }
}

Remove-HyperlinkAtAddress $ws '$A$6'
Remove-HyperlinkAtAddress $ws '$A$7'

# --- Re-create a single hyperlink on the (new) A7, which now holds
# "raj@abc.com" after the swap above ---
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:raj@abc.com")

# Restore the original (pre-existing) cell styles so no new style gets
# left applied on the cells themselves.
$ws.Range("A6").Style = $styleA6
$ws.Range("A7").Style = $styleA7

# --- Update the selection to match: whole row 6 selected, active cell A6 ---
$ws.Rows.Item(6).Select()
